$wb = $excel.ActiveWorkbook

# "想去人数" (want-to-go headcount) counts were refreshed for the data export.
# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 469
$wsExhibition.Range("F4").Value = 7854
$wsExhibition.Range("F5").Value = 94
$wsExhibition.Range("F6").Value = 215
$wsExhibition.Range("F10").Value = 457
$wsExhibition.Range("F11").Value = 165
$wsExhibition.Range("F13").Value = 444
$wsExhibition.Range("F14").Value = 66
$wsExhibition.Range("F15").Value = 69
$wsExhibition.Range("F17").Value = 5747
$wsExhibition.Range("F18").Value = 167
$wsExhibition.Range("F20").Value = 1462
$wsExhibition.Range("F22").Value = 354

# Sheet "全部类型" (All types) mirrors the same events at different row offsets.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 469
$wsAll.Range("F4").Value = 7854
$wsAll.Range("F5").Value = 94
$wsAll.Range("F6").Value = 215
$wsAll.Range("F10").Value = 457
$wsAll.Range("F11").Value = 165
$wsAll.Range("F13").Value = 444
$wsAll.Range("F14").Value = 66
$wsAll.Range("F15").Value = 69
$wsAll.Range("F18").Value = 5747
$wsAll.Range("F20").Value = 167
$wsAll.Range("F22").Value = 1462
$wsAll.Range("F24").Value = 354
